# Commit: "1st changes of mifos to finflux"
# The "Repayment Schedule" sheet gains a new blank column inserted right
# before column N ("Late"), shifting the old N:P block (Late / blank / Outstanding)
# one column to the right (becoming O:Q). Everything else on the sheet is
# unaffected by the insert - styles/widths/row data all shift along with it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert one blank column immediately before column N; this pushes the
# existing N, O, P columns (and all 14 rows of data in them) one column
# to the right, to O, P, Q respectively - matching the diff exactly.
$ws.Columns("N:N").Insert()

# Reflect the new selection/cursor position recorded in the saved file.
$ws.Range("S6").Select()
